# Auto-generated script to apply scheduled market-price refresh updates
# across the Jenova Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 993.36664
$ws.Range("I15").Value = 993.36664
$ws.Range("K15").Value = 2980.09992
$ws.Range("M15").Value = -2811.09992
$ws.Range("H74").Value = 14649.5
$ws.Range("I74").Value = 4299
$ws.Range("K74").Value = 4299
$ws.Range("M74").Value = -3363
$ws.Range("H77").Value = 14649.5
$ws.Range("I77").Value = 4299
$ws.Range("K77").Value = 21495
$ws.Range("M77").Value = -16815
$ws.Range("H92").Value = 360.11765
$ws.Range("I92").Value = 316.41666
$ws.Range("J92").Value = 465
$ws.Range("K92").Value = 316.41666
$ws.Range("L92").Value = 465
$ws.Range("M92").Value = 931.58334
$ws.Range("N92").Value = -2961
$ws.Range("H96").Value = 907.4
$ws.Range("I96").Value = 911
$ws.Range("K96").Value = 2733
$ws.Range("M96").Value = -1360
$ws.Range("H99").Value = 84527.664
$ws.Range("J99").Value = 168855.5
$ws.Range("L99").Value = 506566.5
$ws.Range("N99").Value = -509562.5
$ws.Range("H101").Value = 1002.4
$ws.Range("I101").Value = 1054.3334
$ws.Range("J101").Value = 924.5
$ws.Range("K101").Value = 3163.0002
$ws.Range("L101").Value = 2773.5
$ws.Range("M101").Value = -1541.0002
$ws.Range("N101").Value = -6017.5
$ws.Range("H104").Value = 2789.5
$ws.Range("I104").Value = 2789.5
$ws.Range("K104").Value = 8368.5
$ws.Range("M104").Value = -6621.5
$ws.Range("H106").Value = 3144.1292
$ws.Range("I106").Value = 3036.08
$ws.Range("J106").Value = 3594.3333
$ws.Range("K106").Value = 3036.08
$ws.Range("L106").Value = 3594.3333
$ws.Range("M106").Value = -2405.08
$ws.Range("N106").Value = -4856.3333
$ws.Range("H112").Value = 1877.6111
$ws.Range("J112").Value = 1877.6111
$ws.Range("L112").Value = 5632.8333
$ws.Range("N112").Value = -7848.8333
$ws.Range("H132").Value = 1598.7931
$ws.Range("I132").Value = 1385.6792
$ws.Range("J132").Value = 3857.8
$ws.Range("K132").Value = 4157.0376
$ws.Range("L132").Value = 11573.4
$ws.Range("M132").Value = -1627.0376
$ws.Range("N132").Value = -16633.4
$ws.Range("H138").Value = 6676.058
$ws.Range("I138").Value = 4121
$ws.Range("J138").Value = 7450.3184
$ws.Range("K138").Value = 12363
$ws.Range("L138").Value = 22350.9552
$ws.Range("M138").Value = -7223
$ws.Range("N138").Value = -32630.9552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2342.7917
$ws.Range("I102").Value = 2275.7273
$ws.Range("K102").Value = 2275.7273
$ws.Range("M102").Value = -653.7273
$ws.Range("H132").Value = 3933.3823
$ws.Range("I132").Value = 3314.5
$ws.Range("J132").Value = 6320.5
$ws.Range("K132").Value = 9943.5
$ws.Range("L132").Value = 18961.5
$ws.Range("M132").Value = -7413.5
$ws.Range("N132").Value = -24021.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2231.1428
$ws.Range("I99").Value = 2406
$ws.Range("J99").Value = 1998
$ws.Range("K99").Value = 2406
$ws.Range("L99").Value = 1998
$ws.Range("M99").Value = -908
$ws.Range("N99").Value = -4994
$ws.Range("H105").Value = 113463.22
$ws.Range("I105").Value = 144879.86
$ws.Range("K105").Value = 144879.86
$ws.Range("M105").Value = -143132.86
$ws.Range("H132").Value = 48259.89
$ws.Range("J132").Value = 48259.89
$ws.Range("L132").Value = 48259.89
$ws.Range("N132").Value = -58379.89
$ws.Range("H134").Value = 21653.637
$ws.Range("I134").Value = 3395.283
$ws.Range("K134").Value = 10185.849
$ws.Range("M134").Value = -7650.849

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32219
$ws.Range("I31").Value = 1740.8096
$ws.Range("J31").Value = 81453
$ws.Range("K31").Value = 1740.8096
$ws.Range("L31").Value = 81453
$ws.Range("M31").Value = -1445.8096
$ws.Range("N31").Value = -82043
$ws.Range("H34").Value = 32219
$ws.Range("I34").Value = 1740.8096
$ws.Range("J34").Value = 81453
$ws.Range("K34").Value = 1740.8096
$ws.Range("L34").Value = 81453
$ws.Range("M34").Value = -1538.8096
$ws.Range("N34").Value = -81857
$ws.Range("H58").Value = 2846.923
$ws.Range("I58").Value = 2805.1177
$ws.Range("J58").Value = 2925.889
$ws.Range("K58").Value = 2805.1177
$ws.Range("L58").Value = 2925.889
$ws.Range("M58").Value = -2602.1177
$ws.Range("N58").Value = -3331.889
$ws.Range("H94").Value = 867.3333
$ws.Range("I94").Value = 286.57144
$ws.Range("K94").Value = 286.57144
$ws.Range("M94").Value = 164.42856
$ws.Range("H105").Value = 998.1818
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 2344.6365
$ws.Range("I122").Value = 966
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 2898
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -448
$ws.Range("N122").Value = -16897
$ws.Range("H132").Value = 4709.6665
$ws.Range("I132").Value = 3064.5
$ws.Range("K132").Value = 9193.5
$ws.Range("M132").Value = -6663.5
$ws.Range("H134").Value = 324762.2
$ws.Range("I134").Value = 2347.2693
$ws.Range("J134").Value = 2001319.8
$ws.Range("K134").Value = 7041.8079
$ws.Range("L134").Value = 6003959.4
$ws.Range("M134").Value = -4506.8079
$ws.Range("N134").Value = -6009029.4
$ws.Range("H136").Value = 2846.923
$ws.Range("I136").Value = 2805.1177
$ws.Range("J136").Value = 2925.889
$ws.Range("K136").Value = 8415.3531
$ws.Range("L136").Value = 8777.667000000001
$ws.Range("M136").Value = -5865.3531
$ws.Range("N136").Value = -13877.667
$ws.Range("H141").Value = 181474.48
$ws.Range("J141").Value = 180177.12
$ws.Range("L141").Value = 180177.12
$ws.Range("N141").Value = -190537.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1350530.2
$ws.Range("J5").Value = 3333785
$ws.Range("L5").Value = 10001355
$ws.Range("N5").Value = -10001579
$ws.Range("H68").Value = 3334167.2
$ws.Range("J68").Value = 3334334
$ws.Range("L68").Value = 10003002
$ws.Range("N68").Value = -10004624
$ws.Range("H71").Value = 3334167.2
$ws.Range("J71").Value = 3334334
$ws.Range("L71").Value = 30009006
$ws.Range("N71").Value = -30017118
$ws.Range("H94").Value = 13900
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H106").Value = 38250
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 120000
$ws.Range("N106").Value = -121892
$ws.Range("H129").Value = 37149704
$ws.Range("J129").Value = 252625
$ws.Range("L129").Value = 757875
$ws.Range("N129").Value = -767875
$ws.Range("H135").Value = 1350530.2
$ws.Range("J135").Value = 3333785
$ws.Range("L135").Value = 30004065
$ws.Range("N135").Value = -30009135

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4348.36
$ws.Range("I102").Value = 2208.923
$ws.Range("J102").Value = 6666.0835
$ws.Range("K102").Value = 2208.923
$ws.Range("L102").Value = 6666.0835
$ws.Range("M102").Value = -586.9229999999998
$ws.Range("N102").Value = -9910.083500000001
$ws.Range("H123").Value = 74992.336
$ws.Range("J123").Value = 74992.336
$ws.Range("L123").Value = 74992.336
$ws.Range("N123").Value = -79892.336
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 106300.6
$ws.Range("I40").Value = 204601.4
$ws.Range("K40").Value = 204601.4
$ws.Range("M40").Value = -204465.4
$ws.Range("H68").Value = 335333.34
$ws.Range("I68").Value = 3000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2251
$ws.Range("H71").Value = 335333.34
$ws.Range("I71").Value = 3000
$ws.Range("K71").Value = 15000
$ws.Range("M71").Value = -11256

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1670.84
$ws.Range("I126").Value = 1599.3889
$ws.Range("K126").Value = 4798.1667
$ws.Range("M126").Value = -2328.1667

Write-Host "Applied scheduled market price updates to all sheets."
